$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Povit")
$ws.Range("A1").Value = "test"
